$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Login with valid username and password"
$ws.Range("B7").Value = "PASSED"
$ws.Range("C7").Value = "chrome"
$ws.Range("D7").Value = "23_01_23212731"

$ws.Range("A8").Value = "Create Country"
$ws.Range("B8").Value = "PASSED"
$ws.Range("C8").Value = "chrome"
$ws.Range("D8").Value = "23_01_23212747"
